$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: column headers for the test case table (filled first to control
# shared-string table ordering)
$ws.Range("A2").Value = "Nmae "

# Row 1: "inputs" header - A1:E1 is already merged, just update its value.
$ws.Range("A1").Value = "inputs"

$ws.Range("B2").Value = "miles traveled"
$ws.Range("C2").Value = "buying artifact"
$ws.Range("D2").Value = "money on them"
$ws.Range("E2").Value = "poster bought"
$ws.Range("F2").Value = "Reward for flyers"

# xlCenter = -4108. Setting alignment cell-by-cell (instead of on the whole
# F1:J1 range at once) keeps the style table from growing duplicate entries.
foreach ($addr in @("F1", "G1", "H1", "I1", "J1")) {
    $ws.Range($addr).HorizontalAlignment = -4108
}
$ws.Range("F1:J1").Merge()

$ws.Range("K1").Value = "outputs"

$ws.Range("J2").Select()
